$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '30.583.46'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.14%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '1.926.25'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +1.33%  '

$ws.Cells.Item(4, 5).Value = '  -0.38%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '246.72'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +2.34%  '

$ws.Cells.Item(6, 5).Value = '  -0.01%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.4738'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -1.28%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.2915'
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -2.11%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.06806'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +2.71%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '105.82'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +4.94%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '18.38'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -3.82%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '1.916.00'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +1.41%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '0.07737'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +1.17%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '5.342'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +3.46%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.6718'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +0.84%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '288.15'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -5.39%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '30.613.35'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.35%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '0.000007653'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.84%  '

$ws.Cells.Item(19, 2).Value = 'Avalanche'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '13.06'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -1.54%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.14%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '2.167.65'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +1.13%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '5.460'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +4.34%  '

$ws.Cells.Item(23, 5).Value = '  +0.10%  '

$ws.Cells.Item(24, 5).Value = '  +0.05%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '9.397'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.02%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '168.76'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.27%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '20.72'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.58%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '2.142'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +6.64%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '0.1088'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -2.43%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '1.360'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +0.28%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '4.171'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -0.42%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '4.025'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +0.25%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '0.05071'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -0.67%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.7407'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -2.03%  '

$ws.Cells.Item(35, 5).Value = '  -1.69%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.02090'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +4.06%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '2.732'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -0.98%  '

$ws.Cells.Item(38, 5).Value = '  -1.26%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '2.062'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.70%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '111.14'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.56%  '

$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '0.4442'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +4.58%  '

$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.8764'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.18%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '5.917'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +2.84%  '

$ws.Cells.Item(44, 5).Value = '  +0.10%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '67.65'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.98%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '7.273'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -1.52%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '9.364'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.56%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '0.1235'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.35%  '

$ws.Cells.Item(49, 2).Value = 'BitcoinSV'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '47.38'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +9.84%  '

$ws.Cells.Item(50, 2).Value = 'Decentraland'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '0.4130'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +6.13%  '

$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '35.20'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.48%  '
